# Applies a set of corrective cell-value updates to the single worksheet in
# the workbook. The underlying data had several pairs (and one triple) of
# adjacent rows whose stock-movement figures (item code, rate, qty, value)
# were transposed; this script rewrites the affected cells with the correct
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 272 / 273 -------------------------------------------------------
$ws.Range("B272").Value = 48706
$ws.Range("E272").Value = 39.8
$ws.Range("F272").Value = -144
$ws.Range("G272").Value = -4795.2

$ws.Range("B273").Value = 64973
$ws.Range("E273").Value = 35.4
$ws.Range("F273").Value = 150
$ws.Range("G273").Value = 4995

# --- rows 306 / 307 -------------------------------------------------------
$ws.Range("B306").Value = 62997
$ws.Range("F306").Value = 72
$ws.Range("G306").Value = 22020.48

$ws.Range("B307").Value = 57854
$ws.Range("F307").Value = 2
$ws.Range("G307").Value = 611.6799999999999

# --- rows 343 / 344 / 345 --------------------------------------------------
$ws.Range("B343").Value = 63531
$ws.Range("E343").Value = 152.53
$ws.Range("F343").Value = 80
$ws.Range("G343").Value = 11478.4

$ws.Range("B344").Value = 57802
$ws.Range("E344").Value = 162.71
$ws.Range("F344").Value = -79
$ws.Range("G344").Value = -11334.92

$ws.Range("B345").Value = 63571
$ws.Range("F345").Value = 29
$ws.Range("G345").Value = 4160.92

# --- rows 375 / 376 -------------------------------------------------------
$ws.Range("B375").Value = 60325
$ws.Range("E375").Value = 151.57
$ws.Range("F375").Value = -102
$ws.Range("G375").Value = -12939.72

$ws.Range("B376").Value = 63560
$ws.Range("E376").Value = 134.87
$ws.Range("F376").Value = 104
$ws.Range("G376").Value = 13193.44

# --- rows 382 / 383 -------------------------------------------------------
$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31

$ws.Range("B383").Value = 57817
$ws.Range("F383").Value = 3
$ws.Range("G383").Value = 239.43

# --- rows 412 / 413 -------------------------------------------------------
$ws.Range("B412").Value = 63007
$ws.Range("F412").Value = 984
$ws.Range("G412").Value = 168588.72

$ws.Range("B413").Value = 57856
$ws.Range("F413").Value = 2
$ws.Range("G413").Value = 342.66

# --- rows 424 / 425 -------------------------------------------------------
$ws.Range("B424").Value = 53082
$ws.Range("C424").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F424").Value = 1
$ws.Range("G424").Value = 59.47

$ws.Range("B425").Value = 63102
$ws.Range("C425").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F425").Value = 36
$ws.Range("G425").Value = 2140.92

# --- rows 450 / 451 -------------------------------------------------------
$ws.Range("B450").Value = 31930
$ws.Range("E450").Value = 26.8
$ws.Range("F450").Value = -62
$ws.Range("G450").Value = -1390.04

$ws.Range("B451").Value = 63681
$ws.Range("E451").Value = 23.84
$ws.Range("F451").Value = 65
$ws.Range("G451").Value = 1457.3

# --- rows 529 / 530 -------------------------------------------------------
$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2

$ws.Range("B530").Value = 58047
$ws.Range("D530").Value = 105.54
$ws.Range("E530").Value = 126.1
$ws.Range("F530").Value = 54
$ws.Range("G530").Value = 5699.16

# --- rows 572 / 573 -------------------------------------------------------
$ws.Range("B572").Value = 53757
$ws.Range("E572").Value = 16.08
$ws.Range("F572").Value = -159
$ws.Range("G572").Value = -2138.55

$ws.Range("B573").Value = 65069
$ws.Range("E573").Value = 14.3
$ws.Range("F573").Value = 172
$ws.Range("G573").Value = 2313.4

# --- rows 579 / 580 -------------------------------------------------------
$ws.Range("B579").Value = 45695
$ws.Range("E579").Value = 23.58
$ws.Range("F579").Value = -36
$ws.Range("G579").Value = -710.28

$ws.Range("B580").Value = 64915
$ws.Range("E580").Value = 20.98
$ws.Range("F580").Value = 40
$ws.Range("G580").Value = 789.2

# --- rows 583 / 584 -------------------------------------------------------
$ws.Range("B583").Value = 45706
$ws.Range("E583").Value = 23.58
$ws.Range("F583").Value = -202
$ws.Range("G583").Value = -3985.46

$ws.Range("B584").Value = 64922
$ws.Range("E584").Value = 20.98
$ws.Range("F584").Value = 207
$ws.Range("G584").Value = 4084.11

# --- rows 586 / 587 -------------------------------------------------------
$ws.Range("B586").Value = 64927
$ws.Range("E586").Value = 17.26
$ws.Range("F586").Value = 295
$ws.Range("G586").Value = 4784.9

$ws.Range("B587").Value = 45718
$ws.Range("E587").Value = 19.38
$ws.Range("F587").Value = -294
$ws.Range("G587").Value = -4768.68

# --- rows 680 / 681 -------------------------------------------------------
$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52

$ws.Range("B681").Value = 64810
$ws.Range("E681").Value = 291.22
$ws.Range("F681").Value = 7
$ws.Range("G681").Value = 1917.44

# --- rows 702 / 703 -------------------------------------------------------
$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34

$ws.Range("B703").Value = 64833
$ws.Range("E703").Value = 34.9
$ws.Range("F703").Value = 99
$ws.Range("G703").Value = 3250.17

# --- rows 713 / 714 -------------------------------------------------------
$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79

$ws.Range("B714").Value = 64830
$ws.Range("E714").Value = 34.9
$ws.Range("F714").Value = 117
$ws.Range("G714").Value = 3841.11

# --- rows 865 / 866 -------------------------------------------------------
$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53

$ws.Range("B866").Value = 65079
$ws.Range("E866").Value = 43.44
$ws.Range("F866").Value = 21
$ws.Range("G866").Value = 858.27
